# Update test.pptx table contents/formatting:
#  - Header "product_category" -> "Product Category" (both tables)
#  - Header "is_discounted" -> "Discounted" (slide 2 table)
#  - "Customer Count" column: "N.0" -> "N" (drop the trailing ".0")
#  - "Avg Price" column: raw float -> "$X.XX" (currency, 2 decimals)
#  - "Discount Rate" column: raw fraction -> "XX.X%" (percentage, 1 decimal)
#  - "Total Quantity" column is left untouched.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - "Aggregation 1" table: product_category | Customer Count |
#           Avg Price | Total Quantity | Discount Rate
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$t1 = $s1.Shapes.Item(2).Table

$t1.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Product Category"

$slide1Rows = @(
    @{ Count = "130"; Price = "`$31.39"; Discount = "28.0%" },
    @{ Count = "226"; Price = "`$34.70"; Discount = "32.3%" },
    @{ Count = "257"; Price = "`$31.50"; Discount = "32.4%" },
    @{ Count = "182"; Price = "`$34.70"; Discount = "28.5%" },
    @{ Count = "112"; Price = "`$35.23"; Discount = "27.7%" }
)

for ($i = 0; $i -lt $slide1Rows.Count; $i++) {
    $row = $i + 2
    $data = $slide1Rows[$i]
    $t1.Cell($row, 2).Shape.TextFrame.TextRange.Text = $data.Count
    $t1.Cell($row, 3).Shape.TextFrame.TextRange.Text = $data.Price
    $t1.Cell($row, 5).Shape.TextFrame.TextRange.Text = $data.Discount
}

# ---------------------------------------------------------------------------
# Slide 2 - "Aggregation 2" table: product_category | is_discounted |
#           Customer Count | Avg Price | Total Quantity | Discount Rate
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$t2 = $s2.Shapes.Item(2).Table

$t2.Cell(1, 1).Shape.TextFrame.TextRange.Text = "Product Category"
$t2.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Discounted"

$slide2Rows = @(
    @{ Count = "94";  Price = "`$30.71"; Discount = "0.0%" },
    @{ Count = "40";  Price = "`$33.13"; Discount = "100.0%" },
    @{ Count = "160"; Price = "`$35.14"; Discount = "0.0%" },
    @{ Count = "81";  Price = "`$33.78"; Discount = "100.0%" },
    @{ Count = "176"; Price = "`$33.32"; Discount = "0.0%" },
    @{ Count = "91";  Price = "`$27.69"; Discount = "100.0%" },
    @{ Count = "129"; Price = "`$35.51"; Discount = "0.0%" },
    @{ Count = "57";  Price = "`$32.66"; Discount = "100.0%" },
    @{ Count = "81";  Price = "`$32.08"; Discount = "0.0%" },
    @{ Count = "31";  Price = "`$43.45"; Discount = "100.0%" }
)

for ($i = 0; $i -lt $slide2Rows.Count; $i++) {
    $row = $i + 2
    $data = $slide2Rows[$i]
    $t2.Cell($row, 3).Shape.TextFrame.TextRange.Text = $data.Count
    $t2.Cell($row, 4).Shape.TextFrame.TextRange.Text = $data.Price
    $t2.Cell($row, 6).Shape.TextFrame.TextRange.Text = $data.Discount
}
